$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1964285714285714
$ws.Range("C2").Value = 0.5848214285714286
$ws.Range("J2").Value = 0.008928571428571428
$ws.Range("P2").Value = 0.1383928571428572
$ws.Range("S2").Value = 0.07142857142857142
$ws.Range("B3").Value = 0.02290076335877863
$ws.Range("J3").Value = 0.01526717557251908
$ws.Range("P3").Value = 0.7404580152671756
$ws.Range("S3").Value = 0.2213740458015267
$ws.Range("P4").Value = 0.7741935483870968
$ws.Range("S4").Value = 0.2258064516129032
$ws.Range("B6").Value = 0.03669724770642202
$ws.Range("D6").Value = 0.009174311926605505
$ws.Range("F6").Value = 0.07798165137614679
$ws.Range("J6").Value = 0.2752293577981652
$ws.Range("O6").Value = 0.004587155963302753
$ws.Range("Q6").Value = 0.1238532110091743
$ws.Range("R6").Value = 0.05504587155963303
$ws.Range("S6").Value = 0.4174311926605505
$ws.Range("B7").Value = 0.06701030927835051
$ws.Range("D7").Value = 0.0154639175257732
$ws.Range("F7").Value = 0.05670103092783505
$ws.Range("J7").Value = 0.1391752577319588
$ws.Range("O7").Value = 0.02061855670103093
$ws.Range("Q7").Value = 0.1907216494845361
$ws.Range("R7").Value = 0.06185567010309279
$ws.Range("S7").Value = 0.4484536082474227
$ws.Range("B8").Value = 0.09705372616984402
$ws.Range("D8").Value = 0.01559792027729636
$ws.Range("E8").Value = 0.001733102253032929
$ws.Range("F8").Value = 0.06065857885615251
$ws.Range("J8").Value = 0.1126516464471404
$ws.Range("O8").Value = 0.01039861351819757
$ws.Range("Q8").Value = 0.1386481802426343
$ws.Range("R8").Value = 0.09358752166377816
$ws.Range("S8").Value = 0.4696707105719237
$ws.Range("B9").Value = 0.06779661016949153
$ws.Range("D9").Value = 0.005649717514124294
$ws.Range("F9").Value = 0.07344632768361582
$ws.Range("J9").Value = 0.1016949152542373
$ws.Range("O9").Value = 0.005649717514124294
$ws.Range("Q9").Value = 0.2146892655367232
$ws.Range("R9").Value = 0.0903954802259887
$ws.Range("S9").Value = 0.4406779661016949
$ws.Range("B10").Value = 0.08980827447023208
$ws.Range("D10").Value = 0.01816347124117054
$ws.Range("F10").Value = 0.07769929364278506
$ws.Range("J10").Value = 0.09788092835519677
$ws.Range("O10").Value = 0.006054490413723511
$ws.Range("Q10").Value = 0.2058526740665994
$ws.Range("R10").Value = 0.09687184661957618
$ws.Range("S10").Value = 0.4076690211907165
$ws.Range("G11").Value = 0.1481481481481481
$ws.Range("J11").Value = 0.0962962962962963
$ws.Range("K11").Value = 0.2148148148148148
$ws.Range("L11").Value = 0.5222222222222223
$ws.Range("S11").Value = 0.01851851851851852
$ws.Range("G12").Value = 0.8041958041958042
$ws.Range("J12").Value = 0.1468531468531468
$ws.Range("K12").Value = 0.006993006993006993
$ws.Range("L12").Value = 0.01398601398601399
$ws.Range("S12").Value = 0.02797202797202797
$ws.Range("G13").Value = 0.7543859649122807
$ws.Range("J13").Value = 0.1929824561403509
$ws.Range("S13").Value = 0.05263157894736842
$ws.Range("F15").Value = 0.01648351648351648
$ws.Range("H15").Value = 0.2472527472527473
$ws.Range("I15").Value = 0.06593406593406594
$ws.Range("J15").Value = 0.3241758241758242
$ws.Range("K15").Value = 0.04945054945054945
$ws.Range("M15").Value = 0.01098901098901099
$ws.Range("O15").Value = 0.05494505494505494
$ws.Range("S15").Value = 0.2307692307692308
$ws.Range("F16").Value = 0.03424657534246575
$ws.Range("H16").Value = 0.2191780821917808
$ws.Range("I16").Value = 0.0821917808219178
$ws.Range("J16").Value = 0.3904109589041096
$ws.Range("K16").Value = 0.07534246575342465
$ws.Range("M16").Value = 0.03424657534246575
$ws.Range("O16").Value = 0.03424657534246575
$ws.Range("S16").Value = 0.1301369863013699
$ws.Range("F17").Value = 0.01538461538461539
$ws.Range("H17").Value = 0.2564102564102564
$ws.Range("I17").Value = 0.1
$ws.Range("J17").Value = 0.3256410256410256
$ws.Range("K17").Value = 0.09743589743589744
$ws.Range("M17").Value = 0.02564102564102564
$ws.Range("O17").Value = 0.07179487179487179
$ws.Range("S17").Value = 0.1076923076923077
$ws.Range("F18").Value = 0.02127659574468085
$ws.Range("H18").Value = 0.2393617021276596
$ws.Range("I18").Value = 0.07446808510638298
$ws.Range("J18").Value = 0.324468085106383
$ws.Range("K18").Value = 0.1170212765957447
$ws.Range("M18").Value = 0.01063829787234043
$ws.Range("O18").Value = 0.07446808510638298
$ws.Range("S18").Value = 0.1382978723404255
$ws.Range("F19").Value = 0.02131018153117601
$ws.Range("H19").Value = 0.281767955801105
$ws.Range("I19").Value = 0.07734806629834254
$ws.Range("J19").Value = 0.2975532754538279
$ws.Range("K19").Value = 0.1049723756906077
$ws.Range("M19").Value = 0.02920284135753749
$ws.Range("O19").Value = 0.06393054459352802
$ws.Range("S19").Value = 0.1239147592738753
